$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 56 ("「凧」" entry) entirely, shifting all subsequent rows up by one.
$ws.Rows.Item(56).Delete()
